$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Coin($row, $coin) {
    $ws.Range("B$row").Value = $coin
}

function Set-Link($row, $link) {
    $ws.Range("C$row").Value = $link
}

function Set-Price($row, $price) {
    # The Price column stores values as plain text in the source data
    # (e.g. "74.744.38", "1.00", "0.0000192"). Force text formatting so
    # Excel's automatic number conversion doesn't strip formatting/zeros,
    # then restore the default "Normal" style so no stray number format
    # is left behind on the cell.
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $price
    $cell.Style = "Normal"
}

function Set-Volume($row, $volume) {
    $ws.Range("E$row").Value = $volume
}

Set-Price  2  "74.908.65"
Set-Volume 2  "  +9.09%  "

Set-Price  3  "2.608.44"
Set-Volume 3  "  +7.44%  "

Set-Volume 4  "  -0.04%  "

Set-Price  5  "187.33"
Set-Volume 5  "  +16.72%  "

Set-Price  6  "583.35"
Set-Volume 6  "  +4.37%  "

Set-Volume 7  "  -0.05%  "

Set-Coin   8  "Dogecoin"
Set-Link   8  "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-Price  8  "0.210"
Set-Volume 8  "  +25.72%  "

Set-Coin   9  "XRP"
Set-Link   9  "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-Price  9  "0.539"
Set-Volume 9  "  +5.80%  "

Set-Price  10 "2.610.10"
Set-Volume 10 "  +7.52%  "

Set-Price  11 "0.163"
Set-Volume 11 "  -0.09%  "

Set-Price  12 "0.365"
Set-Volume 12 "  +10.44%  "

Set-Price  13 "4.82"
Set-Volume 13 "  +4.75%  "

Set-Price  14 "0.0000193"
Set-Volume 14 "  +10.45%  "

Set-Price  15 "74.758.18"
Set-Volume 15 "  +9.10%  "

Set-Price  16 "3.087.46"
Set-Volume 16 "  +7.48%  "

Set-Price  17 "26.44"
Set-Volume 17 "  +14.40%  "

Set-Price  18 "2.609.31"
Set-Volume 18 "  +7.50%  "

Set-Price  19 "9.38"
Set-Volume 19 "  +35.89%  "

Set-Price  20 "11.87"
Set-Volume 20 "  +13.36%  "

Set-Price  21 "380.25"
Set-Volume 21 "  +13.55%  "

Set-Price  22 "2.32"
Set-Volume 22 "  +19.55%  "

Set-Price  23 "4.10"
Set-Volume 23 "  +7.45%  "

Set-Price  24 "0.999"
Set-Volume 24 "  -0.08%  "

Set-Price  25 "70.37"
Set-Volume 25 "  +5.18%  "

Set-Price  26 "4.24"
Set-Volume 26 "  +15.48%  "

Set-Price  27 "9.37"
Set-Volume 27 "  +14.38%  "

Set-Price  28 "2.742.88"
Set-Volume 28 "  +7.35%  "

Set-Price  29 "1.00"
Set-Volume 29 "  +0.10%  "

Set-Price  30 "0.0₃0957"
Set-Volume 30 "  +17.00%  "

Set-Coin   31 "Fetch.AI"
Set-Link   31 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-Price  31 "1.40"
Set-Volume 31 "  +21.68%  "

Set-Coin   32 "InternetComputer(DFINITY)"
Set-Link   32 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-Price  32 "8.00"
Set-Volume 32 "  +12.53%  "

Set-Coin   33 "Bittensor"
Set-Link   33 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-Price  33 "512.02"
Set-Volume 33 "  +19.84%  "

Set-Price  34 "1.76"
Set-Volume 34 "  +8.97%  "

Set-Price  35 "0.999"
Set-Volume 35 "  -0.06%  "

Set-Price  36 "0.121"
Set-Volume 36 "  +14.88%  "

Set-Price  37 "159.12"
Set-Volume 37 "  -0.92%  "

Set-Price  38 "19.32"
Set-Volume 38 "  +7.88%  "

Set-Price  39 "19.39"
Set-Volume 39 "  +1.97%  "

Set-Price  41 "4.96"
Set-Volume 41 "  +14.61%  "

Set-Price  42 "1.71"
Set-Volume 42 "  +14.49%  "

Set-Price  43 "0.327"
Set-Volume 43 "  +10.22%  "

Set-Price  44 "2.47"
Set-Volume 44 "  +20.97%  "

Set-Price  45 "157.11"
Set-Volume 45 "  +19.08%  "

Set-Price  46 "1.18"
Set-Volume 46 "  +10.37%  "

Set-Price  47 "38.85"
Set-Volume 47 "  +4.09%  "

Set-Price  48 "0.0833"
Set-Volume 48 "  +16.37%  "

Set-Price  49 "3.64"
Set-Volume 49 "  +9.12%  "

Set-Price  50 "0.529"
Set-Volume 50 "  +9.92%  "

Set-Price  51 "20.14"
Set-Volume 51 "  +19.60%  "
